$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts existing rows 8-12 down to 9-13)
$ws.Rows.Item(8).Insert()

# Match row 8's height to the other data rows (same pattern as row 7/9/11)
$ws.Rows.Item(8).RowHeight = 24.75

# Copy formatting from row 7 (a fully-styled data row) onto the new row 8
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)

# Recreate the merged cells for row 8 to match the other data rows
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

# Fill in the new row's data (serial number, item name, balances, qty, prices)
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "CONCOR PLUS 5/12.5MG 30 F.C. TABLETS"
$ws.Range("H8").Value = "2:2"
$ws.Range("L8").Value = 1
$ws.Range("N8").Value = "72.00"
$ws.Range("P8").Value = "47.5200"
$ws.Range("Q8").Value = "0:2"

# Update the running total (now on row 12 after the insert) to include the new line
$ws.Range("P12").Value = 184.94

# Update the generated timestamp in the footer (now on row 13 after the insert)
$ws.Range("A13").Value = "Wednesday, 24 September, 2025 11:43 AM"
